# add rest of custom fields that were left out
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# These "custom field" cells in row 2 were previously left blank; fill
# them in with their values now that the rest of the fields are wired up.
$ws.Range("J2").Value = "A test creator"        # Creator
$ws.Range("P2").Value = "Some test relation"    # Relation
$ws.Range("Q2").Value = "A test string"         # Temporal
$ws.Range("R2").Value = "Another test one"      # Spatial
$ws.Range("S2").Value = "Something unrelated"   # Spatial Relation

# Move the view/selection onto the newly-populated area.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J2").Select()
